$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 9; $r++) {
    $v = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 2).Value = $v
}
